$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: repeat the variant-number header (R3:U3 = 2,3,4,5) below the table
$ws.Range("R34").Value = 2
$ws.Range("S34").Value = 3
$ws.Range("T34").Value = 4
$ws.Range("U34").Value = 5

# Row 36: column totals over the data rows (R4:U32)
$ws.Range("R36").Formula = "=SUM(R4:R32)"
$ws.Range("S36:U36").Formula = "=SUM(S4:S32)"

# Leave the selection where the author ended up editing
$ws.Range("X27").Select()
